# Updated cryptos list on Mon May  1 03:22:47 UTC 2023 with GitHub Actions
#
# The Price/Volume columns store plain-text numbers (e.g. "336.17",
# "1.005") as strings, not numeric cells. For values that look like a
# normal decimal number, Excel's Range.Value setter auto-converts the
# string to a real number on assignment, so for those cells we briefly
# force the cell to "Text" format, assign the string, then clear the
# format again (ClearFormats) to land back on the default/general look
# while keeping the stored value as text. Values that can't parse as a
# plain number (e.g. "28.665.84", percentages, coin names, URLs) are
# assigned directly since no such coercion risk exists for them.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.665.84'
$ws.Range('E2').Value = '  -2.75%  '
$ws.Range('D3').Value = '1.854.65'
$ws.Range('E3').Value = '  -3.27%  '
$ws.Range('E4').Value = '  -0.87%  '
$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '336.17'
$c.ClearFormats()
$ws.Range('E5').Value = '  +3.34%  '
$c = $ws.Range('D6')
$c.NumberFormat = '@'
$c.Value = '1.005'
$c.ClearFormats()
$c = $ws.Range('D7')
$c.NumberFormat = '@'
$c.Value = '0.4652'
$c.ClearFormats()
$ws.Range('E7').Value = '  -3.21%  '
$c = $ws.Range('D8')
$c.NumberFormat = '@'
$c.Value = '0.3924'
$c.ClearFormats()
$ws.Range('E8').Value = '  -3.07%  '
$c = $ws.Range('D9')
$c.NumberFormat = '@'
$c.Value = '46.45'
$c.ClearFormats()
$ws.Range('E9').Value = '  -3.00%  '
$c = $ws.Range('D10')
$c.NumberFormat = '@'
$c.Value = '0.07929'
$c.ClearFormats()
$ws.Range('E10').Value = '  -3.33%  '
$c = $ws.Range('D11')
$c.NumberFormat = '@'
$c.Value = '0.9844'
$c.ClearFormats()
$ws.Range('E11').Value = '  -2.29%  '
$c = $ws.Range('D12')
$c.NumberFormat = '@'
$c.Value = '22.12'
$c.ClearFormats()
$ws.Range('E12').Value = '  -5.34%  '
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').Value = '1.899.63'
$ws.Range('E13').Value = '  -1.24%  '
$ws.Range('B14').Value = 'Polkadot'
$ws.Range('C14').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$c = $ws.Range('D14')
$c.NumberFormat = '@'
$c.Value = '5.850'
$c.ClearFormats()
$ws.Range('E14').Value = '  -3.31%  '
$c = $ws.Range('D15')
$c.NumberFormat = '@'
$c.Value = '7.026'
$c.ClearFormats()
$ws.Range('E15').Value = '  -2.67%  '
$c = $ws.Range('D16')
$c.NumberFormat = '@'
$c.Value = '0.06795'
$c.ClearFormats()
$ws.Range('E16').Value = '  -0.61%  '
$c = $ws.Range('D17')
$c.NumberFormat = '@'
$c.Value = '1.006'
$c.ClearFormats()
$ws.Range('E17').Value = '  -0.74%  '
$c = $ws.Range('D18')
$c.NumberFormat = '@'
$c.Value = '87.72'
$c.ClearFormats()
$ws.Range('E18').Value = '  -3.88%  '
$ws.Range('E19').Value = '  -2.29%  '
$c = $ws.Range('D20')
$c.NumberFormat = '@'
$c.Value = '17.06'
$c.ClearFormats()
$ws.Range('E20').Value = '  -2.51%  '
$c = $ws.Range('D21')
$c.NumberFormat = '@'
$c.Value = '1.005'
$c.ClearFormats()
$ws.Range('E21').Value = '  -0.61%  '
$ws.Range('D22').Value = '28.652.43'
$ws.Range('E22').Value = '  -2.81%  '
$c = $ws.Range('D23')
$c.NumberFormat = '@'
$c.Value = '5.417'
$c.ClearFormats()
$ws.Range('E23').Value = '  -4.26%  '
$c = $ws.Range('D24')
$c.NumberFormat = '@'
$c.Value = '11.35'
$c.ClearFormats()
$ws.Range('E24').Value = '  -4.28%  '
$c = $ws.Range('D25')
$c.NumberFormat = '@'
$c.Value = '2.134'
$c.ClearFormats()
$ws.Range('E25').Value = '  -2.70%  '
$ws.Range('D26').Value = '2.094.22'
$ws.Range('E26').Value = '  -3.53%  '
$c = $ws.Range('D27')
$c.NumberFormat = '@'
$c.Value = '153.31'
$c.ClearFormats()
$ws.Range('E27').Value = '  -1.77%  '
$c = $ws.Range('D28')
$c.NumberFormat = '@'
$c.Value = '6.263'
$c.ClearFormats()
$ws.Range('E28').Value = '  -3.37%  '
$c = $ws.Range('D29')
$c.NumberFormat = '@'
$c.Value = '19.49'
$c.ClearFormats()
$ws.Range('E29').Value = '  -2.57%  '
$c = $ws.Range('D30')
$c.NumberFormat = '@'
$c.Value = '2.033'
$c.ClearFormats()
$ws.Range('E30').Value = '  -2.93%  '
$c = $ws.Range('D31')
$c.NumberFormat = '@'
$c.Value = '117.62'
$c.ClearFormats()
$ws.Range('E31').Value = '  -2.34%  '
$c = $ws.Range('D32')
$c.NumberFormat = '@'
$c.Value = '0.9833'
$c.ClearFormats()
$ws.Range('E32').Value = '  -2.65%  '
$c = $ws.Range('D33')
$c.NumberFormat = '@'
$c.Value = '0.09449'
$c.ClearFormats()
$ws.Range('E33').Value = '  -1.67%  '
$c = $ws.Range('D34')
$c.NumberFormat = '@'
$c.Value = '5.400'
$c.ClearFormats()
$ws.Range('E34').Value = '  -3.72%  '
$c = $ws.Range('D35')
$c.NumberFormat = '@'
$c.Value = '3.509'
$c.ClearFormats()
$ws.Range('E35').Value = '  -1.35%  '
$c = $ws.Range('D36')
$c.NumberFormat = '@'
$c.Value = '1.348'
$c.ClearFormats()
$ws.Range('E36').Value = '  -1.50%  '
$c = $ws.Range('D37')
$c.NumberFormat = '@'
$c.Value = '0.06159'
$c.ClearFormats()
$ws.Range('E37').Value = '  -2.35%  '
$c = $ws.Range('D38')
$c.NumberFormat = '@'
$c.Value = '0.02201'
$c.ClearFormats()
$ws.Range('E38').Value = '  -3.48%  '
$c = $ws.Range('D39')
$c.NumberFormat = '@'
$c.Value = '1.158'
$c.ClearFormats()
$ws.Range('E39').Value = '  -2.09%  '
$c = $ws.Range('D40')
$c.NumberFormat = '@'
$c.Value = '0.5735'
$c.ClearFormats()
$ws.Range('E40').Value = '  -3.18%  '
$c = $ws.Range('D41')
$c.NumberFormat = '@'
$c.Value = '7.627'
$c.ClearFormats()
$ws.Range('E41').Value = '  -2.94%  '
$c = $ws.Range('D42')
$c.NumberFormat = '@'
$c.Value = '10.11'
$c.ClearFormats()
$ws.Range('E42').Value = '  -5.50%  '
$c = $ws.Range('D43')
$c.NumberFormat = '@'
$c.Value = '0.1790'
$c.ClearFormats()
$ws.Range('E43').Value = '  -3.02%  '
$ws.Range('B44').Value = 'RenderToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$c = $ws.Range('D44')
$c.NumberFormat = '@'
$c.Value = '2.353'
$c.ClearFormats()
$ws.Range('E44').Value = '  -1.56%  '
$ws.Range('B45').Value = 'WEMIXToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$c = $ws.Range('D45')
$c.NumberFormat = '@'
$c.Value = '1.253'
$c.ClearFormats()
$ws.Range('E45').Value = '  -2.25%  '
$ws.Range('B46').Value = 'EnergySwap'
$ws.Range('C46').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$c = $ws.Range('D46')
$c.NumberFormat = '@'
$c.Value = '11.90'
$c.ClearFormats()
$ws.Range('E46').Value = '  -4.14%  '
$ws.Range('B47').Value = 'Decentraland'
$ws.Range('C47').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$c = $ws.Range('D47')
$c.NumberFormat = '@'
$c.Value = '0.5414'
$c.ClearFormats()
$ws.Range('E47').Value = '  -2.47%  '
$ws.Range('E48').Value = '  -4.36%  '
$c = $ws.Range('D49')
$c.NumberFormat = '@'
$c.Value = '1.913'
$c.ClearFormats()
$ws.Range('E49').Value = '  -0.87%  '
$c = $ws.Range('D50')
$c.NumberFormat = '@'
$c.Value = '115.60'
$c.ClearFormats()
$ws.Range('E50').Value = '  -1.97%  '
$c = $ws.Range('D51')
$c.NumberFormat = '@'
$c.Value = '43.59'
$c.ClearFormats()
$ws.Range('E51').Value = '  +4.12%  '
